# Lesson 78 homework/vocab edits: insert the missing business-adjective answers
# into each fill-in-the-blank sentence (split inside the ellipsis run of dots),
# fix a few pre-existing typos, and merge two paragraphs that become one sentence.

$d = $word.ActiveDocument

# 0: My boss is extremely (typo fix + split)
$result0 = $d.Content.Find.Execute("My boss is extremelly ……", $true, $false, $false, $false, $false, $true, 1, $false, "My boss is extremely ……", 2)
Write-Output "[0] My boss is extremely (typo fix + split): $result0"

# 1: Only a ...relentless... person
$result1 = $d.Content.Find.Execute("Only a …………………………person can achieve true success", $true, $false, $false, $false, $false, $true, 1, $false, "Only a …relentless………………………person can achieve true success", 2)
Write-Output "[1] Only a ...relentless... person: $result1"

# 2: Software ...flawless...
$result2 = $d.Content.Find.Execute("Software that we finally implemented was …………………..There was nothing to pick on", $true, $false, $false, $false, $false, $true, 1, $false, "Software that we finally implemented was …flawless………………..There was nothing to pick on", 2)
Write-Output "[2] Software ...flawless...: $result2"

# 3: 1. They offered ...striking...
$result3 = $d.Content.Find.Execute("1. They offered ………………………..amount of money for our company", $true, $false, $false, $false, $false, $true, 1, $false, "1. They offered …striking……………………..amount of money for our company", 2)
Write-Output "[3] 1. They offered ...striking...: $result3"

# 4: 2. The ...staggering... (merge with next paragraph)
$result4 = $d.Content.Find.Execute("2. The………………………..cause of our disastrous situation barriers `r created my leaders on our market", $true, $false, $false, $false, $false, $true, 1, $false, "2. The…staggering……………………..cause of our disastrous situation barriers created my leaders on our market", 2)
Write-Output "[4] 2. The ...staggering... (merge with next paragraph): $result4"

# 5: . Magazine has published ...compromising...
$result5 = $d.Content.Find.Execute(". Magazine has published …………………………..article about one of the most prominent politicians in Poland", $true, $false, $false, $false, $false, $true, 1, $false, ". Magazine has published …compromising………………………..article about one of the most prominent politicians in Poland", 2)
Write-Output "[5] . Magazine has published ...compromising...: $result5"

# 6: . It was a ...smashing... victory
$result6 = $d.Content.Find.Execute(". It was a ………………………….victory. We won this tender leaving everybody behind", $true, $false, $false, $false, $false, $true, 1, $false, ". It was a …smashing……………………….victory. We won this tender leaving everybody behind", 2)
Write-Output "[6] . It was a ...smashing... victory: $result6"

# 7: 5. Economic crisis ...aggravating... sales
$result7 = $d.Content.Find.Execute("5. Economic crisis has been one of the reasons of ………………………sales. Our sales has plummeted", $true, $false, $false, $false, $false, $true, 1, $false, "5. Economic crisis has been one of the reasons of … aggravating ……………………sales. Our sales has plummeted", 2)
Write-Output "[7] 5. Economic crisis ...aggravating... sales: $result7"

# 8: 6. Our company is ...thriving...
$result8 = $d.Content.Find.Execute("6. Our company is ……………………….at the moment. We have reached turnover of 5,000 0000 zl. Our sales has gone through the roof and surpassed our expectations", $true, $false, $false, $false, $false, $true, 1, $false, "6. Our company is … thriving …………………….at the moment. We have reached turnover of 5,000 0000 zl. Our sales has gone through the roof and surpassed our expectations", 2)
Write-Output "[8] 6. Our company is ...thriving...: $result8"

# 9: 7. There is a ...underlying... difference
$result9 = $d.Content.Find.Execute("7.There is a …………………………..difference between our approach to business ethics and approach of our competitors", $true, $false, $false, $false, $false, $true, 1, $false, "7.There is a …underlying ………………………..difference between our approach to business ethics and approach of our competitors", 2)
Write-Output "[9] 7. There is a ...underlying... difference: $result9"

# 10: 8. Their decision ... still ...puzzling... to me
$result10 = $d.Content.Find.Execute("8. Their decision to grant this project to such inexperienced company is still…………………………..to me", $true, $false, $false, $false, $false, $true, 1, $false, "8. Their decision to grant this project to such inexperienced company is still……puzzling……………………..to me", 2)
Write-Output "[10] 8. Their decision ... still ...puzzling... to me: $result10"

# 11: Their arguments occurred to be ...implausible... + didn't typo fix
$result11 = $d.Content.Find.Execute("Their arguments occurred to be ……………………………..They simply did’t convince us", $true, $false, $false, $false, $false, $true, 1, $false, "Their arguments occurred to be …… implausible ………………………..They simply didn’t convince us", 2)
Write-Output "[11] Their arguments occurred to be ...implausible... + didn't typo fix: $result11"

# 12: ew occur(r)ed to be ...reconcilable...
$result12 = $d.Content.Find.Execute("ew occured to be……………………..", $true, $false, $false, $false, $false, $true, 1, $false, "ew occurred to be… reconcilable …………………..", 2)
Write-Output "[12] ew occur(r)ed to be ...reconcilable...: $result12"

# 13: She has ...inexhaustible... energy
$result13 = $d.Content.Find.Execute("She has …………………………energy . She just can’t stop working", $true, $false, $false, $false, $false, $true, 1, $false, "She has …… inexhaustible ……………………energy . She just can’t stop working", 2)
Write-Output "[13] She has ...inexhaustible... energy: $result13"

# 14: Some trade fairs are simply ...unmissable...
$result14 = $d.Content.Find.Execute("Some trade fairs are simply ……………………………You just cannot miss it", $true, $false, $false, $false, $false, $true, 1, $false, "Some trade fairs are simply …… unmissable ………………………You just cannot miss it", 2)
Write-Output "[14] Some trade fairs are simply ...unmissable...: $result14"

# 15: Although this project ... It is ...feasible...
$result15 = $d.Content.Find.Execute("Although this project seems sophisticated It is …………………….meaning that we are capable of executing it", $true, $false, $false, $false, $false, $true, 1, $false, "Although this project seems sophisticated It is … feasible ………………….meaning that we are capable of executing it", 2)
Write-Output "[15] Although this project ... It is ...feasible...: $result15"

# 16: Despite the fact ... remained absolutely ...unaffected...
$result16 = $d.Content.Find.Execute("Despite the fact that this phase was a complete failure he remained absolutely ……………………………..", $true, $false, $false, $false, $false, $true, 1, $false, "Despite the fact that this phase was a complete failure he remained absolutely … unaffected …………………………..", 2)
Write-Output "[16] Despite the fact ... remained absolutely ...unaffected...: $result16"

# 17: When my company went bankrupt ... became ...shattered... (+ typo fixes)
$result17 = $d.Content.Find.Execute("When my company went bunkrupt and bailiff took my possesions I became ……………………….", $true, $false, $false, $false, $false, $true, 1, $false, "When my company went bankrupt and bailiff took my possessions I became … shattered …………………….", 2)
Write-Output "[17] When my company went bankrupt ... became ...shattered... (+ typo fixes): $result17"

# 18: During the last tender ... witness of ...unprecedented... situation
$result18 = $d.Content.Find.Execute("During the last tender I was a witness of ………………………..situation. Foreign company was chosen to conduct one of the largest project in our country", $true, $false, $false, $false, $false, $true, 1, $false, "During the last tender I was a witness of ……unprecedented…………………..situation. Foreign company was chosen to conduct one of the largest project in our country", 2)
Write-Output "[18] During the last tender ... witness of ...unprecedented... situation: $result18"

